# edit.ps1
# Applies the diff: adds Word-style proofing marks (w:proofErr) and splits runs
# accordingly for paragraphs 1, 3, 4, 5; and appends a brand-new paragraph
# (13.01.23 ...) right after paragraph 5, before the trailing empty paragraph.

$d = $word.ActiveDocument

function Replace-ParagraphContent($paraIndex, $xmlFragment) {
    $p = $d.Paragraphs.Item($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the paragraph mark
    $r = $d.Range($start, $end)
    $r.InsertXML($xmlFragment)
}

# --- Paragraph 1: "1.12.22 Отдали старый D-link из кабинета Вики ... wi-fi. ..."
$xmlPara1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>1.12.22</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve"> О</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">тдали старый </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>D</w:t></w:r><w:r w:rsidRPr="007C76A8"><w:t>-</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>link</w:t></w:r><w:r w:rsidRPr="007C76A8"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">из кабинета Вики (лежал на шкафу </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>медниковых</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 9) в 10 кабинет на раздачу </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>wi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="007C76A8"><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>fi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="007C76A8"><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Старый роутер из 10 кабинета оставили там же</w:t></w:r></w:p>'
Replace-ParagraphContent 1 $xmlPara1

# --- Paragraph 3: "6.12.22 На новое рабочее место ..."
$xmlPara3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>6.12.22</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve"> Н</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">а новое рабочее место (где раньше сидел Миша из архитектуры) перенесли </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>комп</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> из 23 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>кбинета</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> и отдали один новый черный сетевой фильтр и свой сетевой кабель.</w:t></w:r></w:p>'
Replace-ParagraphContent 3 $xmlPara3

# --- Paragraph 4: "8.12.22 Отдал свою новую мышку ..."
$xmlPara4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>8.12.22</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve"> О</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">тдал свою новую мышку </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Бархатовой</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Ю.М.  кабинет №12 (у Димы мышка уже была записана на </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Мясникову</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, хотя у нее старая мышь) Мышь </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Бархатовой</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> пока у нас на запасном стенде</w:t></w:r></w:p>'
Replace-ParagraphContent 4 $xmlPara4

# --- Paragraph 5: "8.12.22 Вернул Бархатовой Ю.М. ..."
$xmlPara5 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>8.12.22</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve"> В</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">ернул </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Бархатовой</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Ю.М. кабинет №12 ее старую мышку, забрал новую. Теперь мы должны в отдел </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Мясниковой</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> одну новую мышку </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ExeGate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00921681"><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>потому что по документам у Димы она уже на них числится</w:t></w:r></w:p>'
Replace-ParagraphContent 5 $xmlPara5

# --- New paragraph: "13.01.23 Заменил системник(на новый ...) ..."
# Insert right after paragraph 5's text content, but before paragraph 5's own
# paragraph mark. This pushes in a brand-new paragraph mark/paragraph and
# leaves the trailing empty paragraph untouched and still last in the body.
$p5 = $d.Paragraphs.Item(5)
$insertPos = $p5.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$xmlPara6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">13.01.23 Заменил </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>системни</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>к</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">на новый из тех 5 у которых диск </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>D</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>на 2</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Tb</w:t></w:r><w:r><w:t>) у Кобяковой Оксаны Игоревны (упр. образованием) – старый убрали в кладовку у них там рядом.</w:t></w:r></w:p>'
$insertRange.InsertXML($xmlPara6)
